$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.210.32'
$ws.Range("E2").Value = '  +1.20%  '

$ws.Range("D3").Value = '3.158.27'
$ws.Range("E3").Value = '  -0.96%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.24'
$ws.Range("E5").Value = '  +2.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.17'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.157.44'
$ws.Range("E8").Value = '  -0.97%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  +1.70%  '

$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.37'
$ws.Range("E11").Value = '  +1.96%  '

$ws.Range("E12").Value = '  +0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("E13").Value = '  +1.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.75'
$ws.Range("E14").Value = '  +4.52%  '

$ws.Range("D15").Value = '3.680.03'
$ws.Range("E15").Value = '  -0.87%  '

$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("D17").Value = '3.165.53'
$ws.Range("E17").Value = '  -0.66%  '

$ws.Range("D18").Value = '63.253.07'
$ws.Range("E18").Value = '  +1.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.55'
$ws.Range("E19").Value = '  -0.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '460.40'
$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.93'
$ws.Range("E21").Value = '  +0.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.695'
$ws.Range("E22").Value = '  -1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.62'
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.24'
$ws.Range("E24").Value = '  -1.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.04'
$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.07'
$ws.Range("E29").Value = '  +2.82%  '

$ws.Range("E30").Value = '  -1.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.72'
$ws.Range("E31").Value = '  -2.75%  '

$ws.Range("E32").Value = '  -1.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0999'
$ws.Range("E33").Value = '  -2.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.40'
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("E35").Value = '  -2.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.88'
$ws.Range("E36").Value = '  +1.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.16'
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").Value = '0.0₃0724'
$ws.Range("E38").Value = '  +4.69%  '

$ws.Range("E39").Value = '  +0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.11'
$ws.Range("E40").Value = '  +1.25%  '

$ws.Range("E41").Value = '  -0.66%  '

$ws.Range("E42").Value = '  -0.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '388.90'
$ws.Range("E43").Value = '  -5.51%  '

$ws.Range("D44").Value = '2.783.17'
$ws.Range("E44").Value = '  -5.31%  '

$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.77'
$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.10'
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.54'
$ws.Range("E49").Value = '  +1.05%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.111'
$ws.Range("E50").Value = '  +0.46%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.95'
$ws.Range("E51").Value = '  -2.13%  '
